$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing totals in row 7 (Tổng Hạng / Tổng Điểm) for HẠNG 2 ĐỘI: GLI block
$ws.Range("E7").Value = 5.5
$ws.Range("G7").Value = 9.5

# Add new team block (HẠNG 9 ĐỘI: TXQ) in rows 27-29
$ws.Range("A27").Value = "HẠNG 9 ĐỘI: TXQ"

$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Nguyễn Ngọc Đăng Khoa"
$ws.Range("C28").Value = "TXQ"
$ws.Range("D28").Value = 9
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = 37
$ws.Range("G28").Value = 8

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Nguyễn Ngọc Đăng Khoa"
$ws.Range("C29").Value = "TXQ"
$ws.Range("D29").Value = 28
$ws.Range("E29").Value = 4
